$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition): update "想去人数" (want-to-go count) values
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 310
$wsExhibition.Range("F4").Value = 1251
$wsExhibition.Range("F5").Value = 626

# Sheet "全部类型" (All Types): same underlying events, different row layout
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 310
$wsAll.Range("F4").Value = 1251
$wsAll.Range("F6").Value = 626
